# fall 23 week 14 inputs
# Appends 47 new matchup rows (rows 2462-2508) to Sheet1, columns A:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(5,0,4,2),
    @(5,2,5,0),
    @(4,0,4,2),
    @(3,3,2,0),
    @(5,3,2,0),
    @(5,0,5,3),
    @(4,3,4,0),
    @(3,2,3,1),
    @(5,3,3,0),
    @(6,2,5,0),
    @(3,2,2,1),
    @(6,2,5,0),
    @(5,2,5,0),
    @(4,2,5,0),
    @(3,2,4,1),
    @(4,2,3,1),
    @(4,0,6,3),
    @(4,3,3,0),
    @(5,0,6,2),
    @(4,0,3,2),
    @(4,1,3,2),
    @(6,1,7,2),
    @(4,1,5,2),
    @(6,2,5,0),
    @(5,3,3,0),
    @(4,2,3,1),
    @(3,1,3,2),
    @(6,0,6,2),
    @(2,0,6,3),
    @(3,2,5,0),
    @(5,2,4,0),
    @(6,3,5,0),
    @(5,2,4,0),
    @(4,0,4,2),
    @(6,0,5,2),
    @(3,3,4,0),
    @(6,3,3,0),
    @(3,2,3,1),
    @(4,2,5,0),
    @(2,3,2,0),
    @(5,0,7,2),
    @(3,2,4,0),
    @(2,2,3,1),
    @(6,3,5,0),
    @(6,3,6,0),
    @(4,0,4,3),
    @(2,1,3,2)
)

$startRow = 2462
$rowIndex = 0
foreach ($record in $data) {
    $r = $startRow + $rowIndex
    $ws.Cells.Item($r, 1).Value = $record[0]
    $ws.Cells.Item($r, 2).Value = $record[1]
    $ws.Cells.Item($r, 3).Value = $record[2]
    $ws.Cells.Item($r, 4).Value = $record[3]
    $rowIndex++
}

$lastRow = $startRow + $data.Count - 1
$nextRow = $lastRow + 1

$ws.Application.ActiveWindow.ScrollRow = $lastRow - 20
$selCell = $ws.Cells.Item($nextRow, 1)
$selCell.Select() | Out-Null
